$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 606.2
$ws.Range("I18").Value = 607.7692
$ws.Range("J18").Value = 596
$ws.Range("K18").Value = 607.7692
$ws.Range("L18").Value = 596
$ws.Range("M18").Value = -323.7692
$ws.Range("N18").Value = -1164
$ws.Range("H19").Value = 1298.4
$ws.Range("I19").Value = 750
$ws.Range("K19").Value = 750
$ws.Range("M19").Value = -575
$ws.Range("H76").Value = 3249.88
$ws.Range("I76").Value = 3036.8047
$ws.Range("J76").Value = 4675.846
$ws.Range("K76").Value = 3036.8047
$ws.Range("L76").Value = 4675.846
$ws.Range("M76").Value = -2721.8047
$ws.Range("N76").Value = -5305.846
$ws.Range("H79").Value = 3249.88
$ws.Range("I79").Value = 3036.8047
$ws.Range("J79").Value = 4675.846
$ws.Range("K79").Value = 3036.8047
$ws.Range("L79").Value = 4675.846
$ws.Range("M79").Value = -1944.8047
$ws.Range("N79").Value = -6859.846
$ws.Range("H113").Value = 12502198
$ws.Range("J113").Value = 25001756
$ws.Range("L113").Value = 25001756
$ws.Range("N113").Value = -25008264
$ws.Range("H129").Value = 963.2353000000001
$ws.Range("J129").Value = 972.6774
$ws.Range("L129").Value = 2918.0322
$ws.Range("N129").Value = -12918.0322
$ws.Range("H132").Value = 1318.74
$ws.Range("I132").Value = 1074.8857
$ws.Range("J132").Value = 1887.7333
$ws.Range("K132").Value = 3224.6571
$ws.Range("L132").Value = 5663.199900000001
$ws.Range("M132").Value = -694.6571000000004
$ws.Range("N132").Value = -10723.1999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 986.0714
$ws.Range("I2").Value = 1066
$ws.Range("J2").Value = 506.5
$ws.Range("K2").Value = 1066
$ws.Range("L2").Value = 506.5
$ws.Range("M2").Value = -953
$ws.Range("N2").Value = -732.5
$ws.Range("H32").Value = 744.09
$ws.Range("I32").Value = 607.25555
$ws.Range("J32").Value = 1975.6
$ws.Range("K32").Value = 607.25555
$ws.Range("L32").Value = 1975.6
$ws.Range("M32").Value = -320.25555
$ws.Range("N32").Value = -2549.6
$ws.Range("H45").Value = 8745.23
$ws.Range("I45").Value = 9390.666999999999
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 9390.666999999999
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -9013.666999999999
$ws.Range("N45").Value = -1754
$ws.Range("H48").Value = 150000
$ws.Range("J48").Value = 150000
$ws.Range("L48").Value = 150000
$ws.Range("N48").Value = -150768
$ws.Range("H61").Value = 2126
$ws.Range("I61").Value = 2346.5454
$ws.Range("K61").Value = 2346.5454
$ws.Range("M61").Value = -2134.5454
$ws.Range("H116").Value = 986.0714
$ws.Range("I116").Value = 1066
$ws.Range("J116").Value = 506.5
$ws.Range("K116").Value = 1066
$ws.Range("L116").Value = 506.5
$ws.Range("M116").Value = 1228
$ws.Range("N116").Value = -5094.5
$ws.Range("H132").Value = 2328517.5
$ws.Range("I132").Value = 1894.9062
$ws.Range("J132").Value = 9096874
$ws.Range("K132").Value = 5684.7186
$ws.Range("L132").Value = 27290622
$ws.Range("M132").Value = -3154.7186
$ws.Range("N132").Value = -27295682
$ws.Range("H136").Value = 2126
$ws.Range("I136").Value = 2346.5454
$ws.Range("K136").Value = 7039.6362
$ws.Range("M136").Value = -4489.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 986.0714
$ws.Range("I3").Value = 1066
$ws.Range("J3").Value = 506.5
$ws.Range("K3").Value = 1066
$ws.Range("L3").Value = 506.5
$ws.Range("M3").Value = -952
$ws.Range("N3").Value = -734.5
$ws.Range("H20").Value = 13168.728
$ws.Range("I20").Value = 1815.9231
$ws.Range("J20").Value = 29567.223
$ws.Range("K20").Value = 1815.9231
$ws.Range("L20").Value = 29567.223
$ws.Range("M20").Value = -1568.9231
$ws.Range("N20").Value = -30061.223
$ws.Range("H99").Value = 47620216
$ws.Range("I99").Value = 62500880
$ws.Range("J99").Value = 2079.8
$ws.Range("K99").Value = 62500880
$ws.Range("L99").Value = 2079.8
$ws.Range("M99").Value = -62499382
$ws.Range("N99").Value = -5075.8
$ws.Range("H134").Value = 1737.0682
$ws.Range("I134").Value = 1125.1786
$ws.Range("J134").Value = 2807.875
$ws.Range("K134").Value = 3375.5358
$ws.Range("L134").Value = 8423.625
$ws.Range("M134").Value = -840.5357999999997
$ws.Range("N134").Value = -13493.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 76668.336
$ws.Range("I4").Value = 90001
$ws.Range("K4").Value = 90001
$ws.Range("M4").Value = -89889
$ws.Range("H16").Value = 1684.8
$ws.Range("I16").Value = 1678.4286
$ws.Range("J16").Value = 1690.375
$ws.Range("K16").Value = 1678.4286
$ws.Range("L16").Value = 1690.375
$ws.Range("M16").Value = -1391.4286
$ws.Range("N16").Value = -2264.375
$ws.Range("H31").Value = 4938
$ws.Range("I31").Value = 1192.3508
$ws.Range("J31").Value = 19171.467
$ws.Range("K31").Value = 1192.3508
$ws.Range("L31").Value = 19171.467
$ws.Range("M31").Value = -897.3507999999999
$ws.Range("N31").Value = -19761.467
$ws.Range("H34").Value = 4938
$ws.Range("I34").Value = 1192.3508
$ws.Range("J34").Value = 19171.467
$ws.Range("K34").Value = 1192.3508
$ws.Range("L34").Value = 19171.467
$ws.Range("M34").Value = -990.3507999999999
$ws.Range("N34").Value = -19575.467
$ws.Range("H113").Value = 1684.8
$ws.Range("I113").Value = 1678.4286
$ws.Range("J113").Value = 1690.375
$ws.Range("K113").Value = 1678.4286
$ws.Range("L113").Value = 1690.375
$ws.Range("M113").Value = 491.5714
$ws.Range("N113").Value = -6030.375
$ws.Range("H122").Value = 1530.1177
$ws.Range("I122").Value = 1044.8889
$ws.Range("J122").Value = 2076
$ws.Range("K122").Value = 3134.6667
$ws.Range("L122").Value = 6228
$ws.Range("M122").Value = -684.6666999999998
$ws.Range("N122").Value = -11128
$ws.Range("H132").Value = 2292
$ws.Range("I132").Value = 1800.2812
$ws.Range("K132").Value = 5400.8436
$ws.Range("M132").Value = -2870.8436
$ws.Range("H134").Value = 1691.2
$ws.Range("I134").Value = 1918.5454
$ws.Range("K134").Value = 5755.6362
$ws.Range("M134").Value = -3220.6362
$ws.Range("H140").Value = 26104.453
$ws.Range("J140").Value = 26104.453
$ws.Range("L140").Value = 26104.453
$ws.Range("N140").Value = -36464.453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 106.77778
$ws.Range("I14").Value = 106.77778
$ws.Range("K14").Value = 320.33334
$ws.Range("M14").Value = -147.33334
$ws.Range("H131").Value = 913.5161000000001
$ws.Range("I131").Value = 698.5714
$ws.Range("J131").Value = 1090.5294
$ws.Range("K131").Value = 2095.7142
$ws.Range("L131").Value = 3271.5882
$ws.Range("M131").Value = 2944.2858
$ws.Range("N131").Value = -13351.5882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 35715244
$ws.Range("I113").Value = 71429270
$ws.Range("J113").Value = 1210
$ws.Range("K113").Value = 71429270
$ws.Range("L113").Value = 1210
$ws.Range("M113").Value = -71427100
$ws.Range("N113").Value = -5550
$ws.Range("H132").Value = 2041.6976
$ws.Range("I132").Value = 1608.12
$ws.Range("J132").Value = 2643.889
$ws.Range("K132").Value = 4824.36
$ws.Range("L132").Value = 7931.667
$ws.Range("M132").Value = -2294.36
$ws.Range("N132").Value = -12991.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 615.55554
$ws.Range("I16").Value = 560
$ws.Range("J16").Value = 726.6667
$ws.Range("K16").Value = 560
$ws.Range("L16").Value = 726.6667
$ws.Range("M16").Value = -390
$ws.Range("N16").Value = -1066.6667
$ws.Range("H26").Value = 8611
$ws.Range("I26").Value = 2222
$ws.Range("J26").Value = 15000
$ws.Range("K26").Value = 2222
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = -1927
$ws.Range("N26").Value = -15590
$ws.Range("H46").Value = 20834154
$ws.Range("I46").Value = 33333878
$ws.Range("J46").Value = 1281.6666
$ws.Range("K46").Value = 33333878
$ws.Range("L46").Value = 1281.6666
$ws.Range("M46").Value = -33333690
$ws.Range("N46").Value = -1657.6666
$ws.Range("H122").Value = 2265991
$ws.Range("I122").Value = 2751899
$ws.Range("J122").Value = 1002629.8
$ws.Range("K122").Value = 8255697
$ws.Range("L122").Value = 3007889.4
$ws.Range("M122").Value = -8253247
$ws.Range("N122").Value = -3012789.4
$ws.Range("H132").Value = 28126612
$ws.Range("I132").Value = 41105830
$ws.Range("J132").Value = 4966.3335
$ws.Range("K132").Value = 123317490
$ws.Range("L132").Value = 14899.0005
$ws.Range("M132").Value = -123314960
$ws.Range("N132").Value = -19959.0005
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H136").Value = 4249.7095
$ws.Range("I136").Value = 1659.0962
$ws.Range("J136").Value = 17720.9
$ws.Range("K136").Value = 4977.2886
$ws.Range("L136").Value = 53162.7
$ws.Range("M136").Value = -2427.2886
$ws.Range("N136").Value = -58262.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1459.8422
$ws.Range("I113").Value = 1373.3846
$ws.Range("J113").Value = 1647.1666
$ws.Range("K113").Value = 4120.1538
$ws.Range("L113").Value = 4941.4998
$ws.Range("M113").Value = -1950.1538
$ws.Range("N113").Value = -9281.4998
$ws.Range("H126").Value = 529.61536
$ws.Range("I126").Value = 365.41666
$ws.Range("K126").Value = 1096.24998
$ws.Range("M126").Value = 1373.75002
$ws.Range("H132").Value = 1039.7593
$ws.Range("I132").Value = 742.4186
$ws.Range("J132").Value = 2202.0908
$ws.Range("K132").Value = 2227.2558
$ws.Range("L132").Value = 6606.2724
$ws.Range("M132").Value = 302.7442000000001
$ws.Range("N132").Value = -11666.2724
$ws.Range("H136").Value = 6758953
$ws.Range("I136").Value = 2299.9363
$ws.Range("J136").Value = 18520534
$ws.Range("K136").Value = 6899.8089
$ws.Range("L136").Value = 55561602
$ws.Range("M136").Value = -4349.8089
$ws.Range("N136").Value = -55566702
